$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H69").Value = 20835828
$ws.Range("I69").Value = 20835828
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 62507484
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -62506610

$ws.Range("H72").Value = 20835828
$ws.Range("I72").Value = 20835828
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 187522452
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -187518084

$ws.Range("H132").Value = 1535.1323
$ws.Range("I132").Value = 1514.3175
$ws.Range("J132").Value = 1797.4
$ws.Range("K132").Value = 4542.9525
$ws.Range("L132").Value = 5392.200000000001
$ws.Range("M132").Value = -2012.9525
$ws.Range("N132").Value = -10452.2

$ws.Range("H137").Value = 4012.1562
$ws.Range("I137").Value = 4936.7856
$ws.Range("J137").Value = 3293
$ws.Range("K137").Value = 14810.3568
$ws.Range("L137").Value = 9879
$ws.Range("M137").Value = -12260.3568
$ws.Range("N137").Value = -14979

$ws.Range("H138").Value = 7154567.5
$ws.Range("I138").Value = 4716
$ws.Range("J138").Value = 20024300
$ws.Range("K138").Value = 14148
$ws.Range("L138").Value = 60072900
$ws.Range("M138").Value = -9008
$ws.Range("N138").Value = -60083180

$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1715993.2
$ws.Range("I32").Value = 1841198.5
$ws.Range("J32").Value = 13201.4
$ws.Range("K32").Value = 1841198.5
$ws.Range("L32").Value = 13201.4
$ws.Range("M32").Value = -1840911.5
$ws.Range("N32").Value = -13775.4

$ws.Range("H61").Value = 6669.1714
$ws.Range("I61").Value = 2645.6191
$ws.Range("J61").Value = 12704.5
$ws.Range("K61").Value = 2645.6191
$ws.Range("L61").Value = 12704.5
$ws.Range("M61").Value = -2433.6191

$ws.Range("H74").Value = 38475.76
$ws.Range("I74").Value = 47342.113
$ws.Range("J74").Value = 10264.637
$ws.Range("K74").Value = 47342.113
$ws.Range("L74").Value = 10264.637
$ws.Range("M74").Value = -46468.113
$ws.Range("N74").Value = -12012.637

$ws.Range("H77").Value = 38475.76
$ws.Range("I77").Value = 47342.113
$ws.Range("J77").Value = 10264.637
$ws.Range("K77").Value = 236710.565
$ws.Range("L77").Value = 51323.185
$ws.Range("M77").Value = -232342.565
$ws.Range("N77").Value = -60059.185

$ws.Range("H136").Value = 6669.1714
$ws.Range("I136").Value = 2645.6191
$ws.Range("J136").Value = 12704.5
$ws.Range("K136").Value = 7936.8573
$ws.Range("L136").Value = 38113.5
$ws.Range("M136").Value = -5386.8573

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 3869.6
$ws.Range("I11").Value = 174.5
$ws.Range("J11").Value = 6333
$ws.Range("K11").Value = 174.5
$ws.Range("L11").Value = 6333
$ws.Range("M11").Value = -34.5
$ws.Range("N11").Value = -6613

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()

$ws.Range("H86").Value = 73591400
$ws.Range("I86").Value = 20918320
$ws.Range("J86").Value = 200006800
$ws.Range("K86").Value = 20918320
$ws.Range("L86").Value = 200006800
$ws.Range("M86").Value = -20917197

$ws.Range("H89").Value = 73591400
$ws.Range("I89").Value = 20918320
$ws.Range("J89").Value = 200006800
$ws.Range("K89").Value = 104591600
$ws.Range("L89").Value = 1000034000
$ws.Range("M89").Value = -104585984

$ws.Range("H94").Value = 3563.6785
$ws.Range("I94").Value = 1854.5
$ws.Range("J94").Value = 5842.5835
$ws.Range("K94").Value = 1854.5
$ws.Range("L94").Value = 5842.5835
$ws.Range("M94").Value = -1403.5
$ws.Range("N94").Value = -6744.5835

$ws.Range("H105").Value = 911727.5600000001
$ws.Range("I105").Value = 1251570
$ws.Range("J105").Value = 5480.8335
$ws.Range("K105").Value = 1251570
$ws.Range("L105").Value = 5480.8335
$ws.Range("M105").Value = -1249823
$ws.Range("N105").Value = -8974.833500000001

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5130.2
$ws.Range("I16").Value = 3133.125
$ws.Range("J16").Value = 7412.5713
$ws.Range("K16").Value = 3133.125
$ws.Range("L16").Value = 7412.5713
$ws.Range("M16").Value = -2846.125
$ws.Range("N16").Value = -7986.5713

$ws.Range("H31").Value = 14126.0625
$ws.Range("I31").Value = 7997.5
$ws.Range("J31").Value = 15001.571
$ws.Range("K31").Value = 7997.5
$ws.Range("L31").Value = 15001.571
$ws.Range("M31").Value = -7702.5
$ws.Range("N31").Value = -15591.571

$ws.Range("H34").Value = 14126.0625
$ws.Range("I34").Value = 7997.5
$ws.Range("J34").Value = 15001.571
$ws.Range("K34").Value = 7997.5
$ws.Range("L34").Value = 15001.571
$ws.Range("M34").Value = -7795.5
$ws.Range("N34").Value = -15405.571

$ws.Range("H113").Value = 5130.2
$ws.Range("I113").Value = 3133.125
$ws.Range("J113").Value = 7412.5713
$ws.Range("K113").Value = 3133.125
$ws.Range("L113").Value = 7412.5713
$ws.Range("M113").Value = -963.125
$ws.Range("N113").Value = -11752.5713

$ws.Range("H132").Value = 5775.727
$ws.Range("I132").Value = 3422.611
$ws.Range("J132").Value = 8599.467000000001
$ws.Range("K132").Value = 10267.833
$ws.Range("L132").Value = 25798.401
$ws.Range("M132").Value = -7737.832999999999
$ws.Range("N132").Value = -30858.401

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 888
$ws.Range("I17").Value = 888
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 2664
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -2495
$ws.Range("N17").ClearContents()

$ws.Range("H68").Value = 4422.0513
$ws.Range("I68").Value = 2871.5
$ws.Range("J68").Value = 4703.9697
$ws.Range("K68").Value = 8614.5
$ws.Range("L68").Value = 14111.9091
$ws.Range("M68").Value = -7803.5
$ws.Range("N68").Value = -15733.9091

$ws.Range("H71").Value = 4422.0513
$ws.Range("I71").Value = 2871.5
$ws.Range("J71").Value = 4703.9697
$ws.Range("K71").Value = 25843.5
$ws.Range("L71").Value = 42335.7273
$ws.Range("M71").Value = -21787.5
$ws.Range("N71").Value = -50447.7273

$ws.Range("H134").Value = 198202.84
$ws.Range("I134").Value = 198202.84
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 594608.52
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -589538.52

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2782.8333
$ws.Range("I80").Value = 3174.5
$ws.Range("J80").Value = 1999.5
$ws.Range("K80").Value = 3174.5
$ws.Range("L80").Value = 1999.5
$ws.Range("M80").Value = -2176.5
$ws.Range("N80").Value = -3995.5

$ws.Range("H83").Value = 2782.8333
$ws.Range("I83").Value = 3174.5
$ws.Range("J83").Value = 1999.5
$ws.Range("K83").Value = 15872.5
$ws.Range("L83").Value = 9997.5
$ws.Range("M83").Value = -10880.5
$ws.Range("N83").Value = -19981.5

$ws.Range("H102").Value = 1695.2858
$ws.Range("I102").Value = 1695.2858
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1695.2858
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -73.28580000000011
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3929.0952
$ws.Range("I93").Value = 3871.5
$ws.Range("J93").Value = 4044.2856
$ws.Range("K93").Value = 3871.5
$ws.Range("L93").Value = 4044.2856
$ws.Range("M93").Value = -2623.5
$ws.Range("N93").Value = -6540.2856

$ws.Range("H107").Value = 4500
$ws.Range("I107").Value = 4500
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4500
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2580

$ws.Range("H122").Value = 2809
$ws.Range("I122").Value = 1912.3334
$ws.Range("J122").Value = 5499
$ws.Range("K122").Value = 5737.0002
$ws.Range("L122").Value = 16497
$ws.Range("M122").Value = -3287.0002
$ws.Range("N122").Value = -21397

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 50015400
$ws.Range("I81").Value = 799.5
$ws.Range("J81").Value = 100030000
$ws.Range("K81").Value = 1599
$ws.Range("L81").Value = 200060000
$ws.Range("M81").Value = -538
$ws.Range("N81").Value = -200062122

$ws.Range("H84").Value = 50015400
$ws.Range("I84").Value = 799.5
$ws.Range("J84").Value = 100030000
$ws.Range("K84").Value = 7995
$ws.Range("L84").Value = 1000300000
$ws.Range("M84").Value = -2691
$ws.Range("N84").Value = -1000310608

$ws.Range("H113").Value = 22111.75
$ws.Range("I113").Value = 36215.715
$ws.Range("J113").Value = 2366.2
$ws.Range("K113").Value = 108647.145
$ws.Range("L113").Value = 7098.599999999999
$ws.Range("M113").Value = -106477.145
$ws.Range("N113").Value = -11438.6

$ws.Range("H122").Value = 8002712
$ws.Range("I122").Value = 10502307
$ws.Range("J122").Value = 4008
$ws.Range("K122").Value = 31506921
$ws.Range("L122").Value = 12024
$ws.Range("M122").Value = -31504471
$ws.Range("N122").Value = -16924

$ws.Range("H132").Value = 13546353
$ws.Range("I132").Value = 16143750
$ws.Range("J132").Value = 126468.164
$ws.Range("K132").Value = 48431250
$ws.Range("L132").Value = 379404.492
$ws.Range("M132").Value = -48428720
$ws.Range("N132").Value = -384464.492

$ws.Range("H136").Value = 58888784
$ws.Range("I136").Value = 200003980
$ws.Range("J136").Value = 90788.5
$ws.Range("K136").Value = 600011940
$ws.Range("L136").Value = 272365.5
$ws.Range("M136").Value = -600009390
$ws.Range("N136").Value = -277465.5
